# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows 470-471) into the Betarraga
# sheet, pushing the existing rows 470-495 down to 472-497.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 470 (shifts old 470:495 -> 472:497)
$ws.Range("A470:A471").EntireRow.Insert()

# New row 470: Primera quality, week of 2023-04-25
$ws.Cells.Item(470, 1).Value = 11
$ws.Cells.Item(470, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(470, 3).Value = "Bíobío"
$ws.Cells.Item(470, 4).Value = 45041
$ws.Cells.Item(470, 5).Value = 8
$ws.Cells.Item(470, 6).Value = 100114014
$ws.Cells.Item(470, 7).Value = "Betarraga"
$ws.Cells.Item(470, 8).Value = "Sin especificar"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 450
$ws.Cells.Item(470, 11).Value = 600
$ws.Cells.Item(470, 12).Value = 650
$ws.Cells.Item(470, 13).Value = 622
$ws.Cells.Item(470, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(470, 15).Value = "Región Metropolitana"
$ws.Cells.Item(470, 16).Value = 124
$ws.Cells.Item(470, 17).Value = 5
$ws.Cells.Item(470, 18).Value = "Hortaliza"

# New row 471: Segunda quality, week of 2023-04-25
$ws.Cells.Item(471, 1).Value = 11
$ws.Cells.Item(471, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(471, 3).Value = "Bíobío"
$ws.Cells.Item(471, 4).Value = 45041
$ws.Cells.Item(471, 5).Value = 8
$ws.Cells.Item(471, 6).Value = 100114014
$ws.Cells.Item(471, 7).Value = "Betarraga"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Segunda"
$ws.Cells.Item(471, 10).Value = 400
$ws.Cells.Item(471, 11).Value = 500
$ws.Cells.Item(471, 12).Value = 500
$ws.Cells.Item(471, 13).Value = 500
$ws.Cells.Item(471, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(471, 15).Value = "Región Metropolitana"
$ws.Cells.Item(471, 16).Value = 100
$ws.Cells.Item(471, 17).Value = 5
$ws.Cells.Item(471, 18).Value = "Hortaliza"
